$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 3.755628166162433

# Row 3
$ws.Range("B3").Value = 0.04172184405617529
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 1.029605918290258

# Row 4
$ws.Range("B4").Value = 1.445647641019636
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 3.223369029078222
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 6.82939032824165

# Row 5
$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.1496068669990043
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 5.582307763322248

# Row 6
$ws.Range("B6").Value = 0.0000119090754144846
$ws.Range("C6").Value = 0.002658071450198252
$ws.Range("D6").Value = 3.223369029078222
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 3.759424968205533
